$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78, shifting existing rows 78..126 down to 79..127.
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new weekly data record.
$ws.Cells.Item(78, 1).Value = 1
$ws.Cells.Item(78, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(78, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(78, 4).Value = 44438
$ws.Cells.Item(78, 5).Value = 15
$ws.Cells.Item(78, 6).Value = "Fruta"
$ws.Cells.Item(78, 7).Value = 100108
$ws.Cells.Item(78, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(78, 9).Value = 100108006
$ws.Cells.Item(78, 10).Value = "Plátano"
$ws.Cells.Item(78, 11).Value = "Sin especificar"
$ws.Cells.Item(78, 12).Value = "Pintón"
$ws.Cells.Item(78, 13).Value = 120
$ws.Cells.Item(78, 14).Value = 19000
$ws.Cells.Item(78, 15).Value = 20000
$ws.Cells.Item(78, 16).Value = 19500
$ws.Cells.Item(78, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(78, 18).Value = "Ecuador"
$ws.Cells.Item(78, 19).Value = 975
$ws.Cells.Item(78, 20).Value = 20
